# "adding position to xlsx"
#
# Inserts a new column D called "position" that concatenates the
# metric (col C) and the group number (col B), e.g. "ASC1_1".
# This pushes the former D:I columns (criteria..text) one column to
# the right, becoming E:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before the current column D (criteria), shifting
# criteria..text from D:I to E:J.
$ws.Columns("D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "position"

# The data lives in three blocks of 9 rows each, separated by a blank
# row (rows 11 and 21 stay blank): 2-10, 12-20, 22-30.
$blockStarts = @(2, 12, 22)
foreach ($startRow in $blockStarts) {
    $endRow = $startRow + 8
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Range("D$r").Formula = '=_xlfn.CONCAT(C' + $r + ',"_",B' + $r + ')'
    }
}

# Reflect the author's final selection (last fill was D22:D30).
$ws.Range("D22:D30").Select() | Out-Null
